# Update data: 2025-10-30 12:44
# Refreshes the "Last Updated" timestamp plus the latest % change / distance
# figures (and associated stock-name reorderings) on the data sheets.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update "Last Updated" timestamp ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Cells.Item(2, 1).Value = "30 Oct 2025, 12:44 PM"

# --- "1 Month Performance" sheet: update stock names and % change values ---
$wsPerf = $wb.Worksheets.Item("1 Month Performance")
$wsPerf.Cells.Item(2, 3).Value = 112.1068
$wsPerf.Cells.Item(4, 3).Value = 78.10429999999999
$wsPerf.Cells.Item(5, 3).Value = 65.52370000000001
$wsPerf.Cells.Item(6, 3).Value = 62.1306
$wsPerf.Cells.Item(7, 2).Value = "MAHASTEEL"
$wsPerf.Cells.Item(7, 3).Value = 53.9642
$wsPerf.Cells.Item(8, 2).Value = "INOXGREEN"
$wsPerf.Cells.Item(8, 3).Value = 53.6377
$wsPerf.Cells.Item(9, 2).Value = "ESSARSHPNG"
$wsPerf.Cells.Item(9, 3).Value = 53.05
$wsPerf.Cells.Item(10, 3).Value = 45.5511
$wsPerf.Cells.Item(11, 3).Value = 41.8479
$wsPerf.Cells.Item(15, 3).Value = 39.0079
$wsPerf.Cells.Item(17, 3).Value = 36.8732
$wsPerf.Cells.Item(20, 3).Value = 34.4939
$wsPerf.Cells.Item(21, 2).Value = "MEGASOFT"
$wsPerf.Cells.Item(21, 3).Value = 33.1719
$wsPerf.Cells.Item(22, 2).Value = "RAMAPHO"
$wsPerf.Cells.Item(22, 3).Value = 33.0381
$wsPerf.Cells.Item(24, 2).Value = "SOUTHBANK"
$wsPerf.Cells.Item(24, 3).Value = 30.2632
$wsPerf.Cells.Item(25, 2).Value = "ORIENTTECH"
$wsPerf.Cells.Item(25, 3).Value = 30.1839
$wsPerf.Cells.Item(26, 2).Value = "MRPL"
$wsPerf.Cells.Item(26, 3).Value = 29.311
$wsPerf.Cells.Item(27, 2).Value = "INDORAMA"
$wsPerf.Cells.Item(27, 3).Value = 29.0804
$wsPerf.Cells.Item(28, 2).Value = "ONMOBILE"
$wsPerf.Cells.Item(28, 3).Value = 28.9142
$wsPerf.Cells.Item(31, 3).Value = 27.2849
$wsPerf.Cells.Item(33, 3).Value = 25.6793
$wsPerf.Cells.Item(34, 3).Value = 25.4813
$wsPerf.Cells.Item(36, 3).Value = 25.2735
$wsPerf.Cells.Item(37, 3).Value = 25.2114
$wsPerf.Cells.Item(38, 3).Value = 24.8156
$wsPerf.Cells.Item(39, 2).Value = "AVALON"
$wsPerf.Cells.Item(39, 3).Value = 24.4138
$wsPerf.Cells.Item(40, 2).Value = "ATHERENERG"
$wsPerf.Cells.Item(40, 3).Value = 24.3034
$wsPerf.Cells.Item(41, 2).Value = "MINDTECK"
$wsPerf.Cells.Item(41, 3).Value = 24.1789
$wsPerf.Cells.Item(42, 3).Value = 24.03
$wsPerf.Cells.Item(43, 2).Value = "UNIPARTS"
$wsPerf.Cells.Item(43, 3).Value = 23.8268
$wsPerf.Cells.Item(44, 2).Value = "AUBANK"
$wsPerf.Cells.Item(44, 3).Value = 23.622
$wsPerf.Cells.Item(45, 3).Value = 23.6168
$wsPerf.Cells.Item(46, 2).Value = "INDIANB"
$wsPerf.Cells.Item(46, 3).Value = 23.336
$wsPerf.Cells.Item(47, 2).Value = "TATVA"
$wsPerf.Cells.Item(47, 3).Value = 22.7182
$wsPerf.Cells.Item(48, 2).Value = "KERNEX"
$wsPerf.Cells.Item(48, 3).Value = 22.6625
$wsPerf.Cells.Item(49, 2).Value = "DCBBANK"
$wsPerf.Cells.Item(49, 3).Value = 22.5591
$wsPerf.Cells.Item(50, 2).Value = "CPEDU"
$wsPerf.Cells.Item(50, 3).Value = 22.3008
$wsPerf.Cells.Item(52, 3).Value = 22.2596
$wsPerf.Cells.Item(53, 3).Value = 22.0921
$wsPerf.Cells.Item(54, 2).Value = "GUJTHEM"
$wsPerf.Cells.Item(54, 3).Value = 21.3115
$wsPerf.Cells.Item(55, 2).Value = "SURYODAY"
$wsPerf.Cells.Item(55, 3).Value = 21.2719
$wsPerf.Cells.Item(56, 3).Value = 21.0714
$wsPerf.Cells.Item(59, 3).Value = 20.3634
$wsPerf.Cells.Item(61, 3).Value = 20.2709
$wsPerf.Cells.Item(62, 3).Value = 20.0898
$wsPerf.Cells.Item(63, 3).Value = 19.9604
$wsPerf.Cells.Item(64, 2).Value = "BANKINDIA"
$wsPerf.Cells.Item(64, 3).Value = 19.646
$wsPerf.Cells.Item(65, 2).Value = "FEDERALBNK"
$wsPerf.Cells.Item(65, 3).Value = 19.6375
$wsPerf.Cells.Item(66, 2).Value = "SHRIRAMFIN"
$wsPerf.Cells.Item(66, 3).Value = 19.5923
$wsPerf.Cells.Item(67, 3).Value = 19.5582
$wsPerf.Cells.Item(68, 2).Value = "PRECWIRE"
$wsPerf.Cells.Item(68, 3).Value = 19.1203
$wsPerf.Cells.Item(69, 2).Value = "MCX"
$wsPerf.Cells.Item(69, 3).Value = 19.0962
$wsPerf.Cells.Item(70, 3).Value = 19.0342
$wsPerf.Cells.Item(71, 2).Value = "AHCL"
$wsPerf.Cells.Item(71, 3).Value = 18.9783
$wsPerf.Cells.Item(73, 2).Value = "WHEELS"
$wsPerf.Cells.Item(73, 3).Value = 18.7279
$wsPerf.Cells.Item(74, 2).Value = "ETHOSLTD"
$wsPerf.Cells.Item(74, 3).Value = 18.6881
$wsPerf.Cells.Item(76, 3).Value = 18.2744

# --- "distance from Dma50" sheet: update distance values ---
$wsDma = $wb.Worksheets.Item("distance from Dma50")
$wsDma.Cells.Item(2, 3).Value = 10.3403
$wsDma.Cells.Item(3, 3).Value = 7.5452
$wsDma.Cells.Item(4, 3).Value = 6.336
$wsDma.Cells.Item(5, 3).Value = 5.3308
$wsDma.Cells.Item(6, 3).Value = 5.2355
$wsDma.Cells.Item(7, 3).Value = 5.0678
$wsDma.Cells.Item(8, 3).Value = 4.5202
$wsDma.Cells.Item(9, 3).Value = 4.393
$wsDma.Cells.Item(10, 3).Value = 3.8676
$wsDma.Cells.Item(11, 3).Value = 3.601
$wsDma.Cells.Item(12, 3).Value = 3.405
$wsDma.Cells.Item(13, 3).Value = 3.4013
$wsDma.Cells.Item(14, 3).Value = 3.1146
$wsDma.Cells.Item(15, 3).Value = 3.0807
$wsDma.Cells.Item(16, 3).Value = 3.003
$wsDma.Cells.Item(17, 3).Value = 2.8549
$wsDma.Cells.Item(18, 3).Value = 2.6715
$wsDma.Cells.Item(19, 3).Value = 2.6476
$wsDma.Cells.Item(20, 3).Value = 2.4149
$wsDma.Cells.Item(21, 3).Value = 2.3663
$wsDma.Cells.Item(22, 3).Value = 1.4805
$wsDma.Cells.Item(23, 3).Value = 1.407
$wsDma.Cells.Item(24, 3).Value = 1.3607
$wsDma.Cells.Item(25, 3).Value = 1.132
$wsDma.Cells.Item(26, 3).Value = 1.0366
$wsDma.Cells.Item(27, 3).Value = 0.9937
$wsDma.Cells.Item(28, 3).Value = 0.6448
$wsDma.Cells.Item(29, 3).Value = 0.2661
$wsDma.Cells.Item(30, 3).Value = -2.0036
